# "Generate Report for Archive"
#
# 1. The localization status for the tracked file moved from
#    "Ready for handoff" to "In Translation". That shared string is
#    surfaced on three sheets:
#      - Overview!E2 and Overview!F2 (per-locale status columns)
#      - zh-cn!C2   (Status column)
#      - de-de!C2   (Status column)
#
# 2. The "Status" column(s) got narrower (report column re-sizing):
#      - Overview columns E ("zh-cn") and F ("de-de")
#      - zh-cn column C ("Status")
#      - de-de column C ("Status")
#    shrank from ~17.22 characters wide down to ~13.41 characters wide.
#    12.5 is the COM ColumnWidth value that lands on the closest
#    achievable pixel-snapped width to the recorded target.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status-related columns ---
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5

Write-Output "Updated status text and resized status columns on Overview, zh-cn, de-de"
